{"js": "const oldTexts = [\n  \"\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u6700\u5173\u952e\u7684\u65b9\u9762\u4e4b\u4e00\u3002 \u5b83\u662f\u4e00\u4e2a\u5e2e\u52a9\u4eba\u4eec\u83b7\u5f97\u77e5\u8bc6\u3001\u6280\u80fd\u3001\u4ef7\u503c\u89c2\u548c\u6001\u5ea6\u7684\u5de5\u5177\uff0c\u8fd9\u4e9b\u5bf9\u4e8e\u4ed6\u4eec\u7684\u4e2a\u4eba\u548c\u804c\u4e1a\u53d1\u5c55\u662f\u5fc5\u8981\u7684\u3002\u6559\u80b2\u5e76\u4e0d\u4ec5\u9650\u4e8e\u5b66\u6821\u548c\u5927\u5b66\uff1b\u5b83\u53ef\u4ee5\u5728\u4efb\u4f55\u5730\u65b9\u548c\u4efb\u4f55\u65f6\u5019\u53d1\u751f\u3002\u6559\u80b2\u662f\u4e00\u4e2a\u6301\u7eed\u7684\u8fc7\u7a0b\uff0c\u5e2e\u52a9\u4e2a\u4eba\u53d1\u6325\u6f5c\u529b\uff0c\u6539\u5584\u751f\u6d3b\u8d28\u91cf\uff0c\u5e76\u4e3a\u793e\u4f1a\u4f5c\u51fa\u8d21\u732e\u3002\",\n  \"\u6559\u80b2\u5bf9\u4e2a\u4eba\u6210\u957f\u548c\u53d1\u5c55\u81f3\u5173\u91cd\u8981\u3002\u5b83\u5e2e\u52a9\u4e2a\u4eba\u66f4\u597d\u5730\u4e86\u89e3\u81ea\u5df1\u7684\u4f18\u70b9\u3001\u7f3a\u70b9\u3001\u5174\u8da3\u548c\u4ef7\u503c\u89c2\uff0c\u4f7f\u4ed6\u4eec\u80fd\u591f\u5236\u5b9a\u76ee\u6807\u3001\u505a\u51fa\u660e\u667a\u51b3\u7b56\u5e76\u5bf9\u81ea\u5df1\u7684\u884c\u52a8\u8d1f\u8d23\u3002\u6559\u80b2\u8fd8\u5e2e\u52a9\u4e2a\u4eba\u57f9\u517b\u6279\u5224\u6027\u601d\u7ef4\u80fd\u529b\u3001\u521b\u9020\u529b\u548c\u89e3\u51b3\u95ee\u9898\u7684\u80fd\u529b\u3002\u8fd9\u4e9b\u6280\u80fd\u5bf9\u4e8e\u5728\u4eca\u5929\u7684\u4e16\u754c\u4e2d\u6210\u529f\u81f3\u5173\u91cd\u8981\uff0c\u8fd9\u4e2a\u4e16\u754c\u4e2d\uff0c\u9002\u5e94\u53d8\u5316\u548c\u89e3\u51b3\u590d\u6742\u95ee\u9898\u7684\u80fd\u529b\u662f\u81f3\u5173\u91cd\u8981\u7684\u3002\",\n  \"\u6559\u80b2\u5728\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u4e2d\u4e5f\u626e\u6f14\u7740\u91cd\u8981\u89d2\u8272\u3002\u5b83\u80fd\u5e2e\u52a9\u4e2a\u4eba\u83b7\u53d6\u5c31\u4e1a\u548c\u521b\u4e1a\u6240\u5fc5\u9700\u7684\u6280\u80fd\u548c\u77e5\u8bc6\u3002\u6559\u80b2\u5bf9\u4e8e\u57f9\u517b\u80fd\u591f\u4e3a\u7ecf\u6d4e\u589e\u957f\u548c\u53d1\u5c55\u8d21\u732e\u7684\u6709\u6280\u80fd\u3001\u6709\u7ade\u4e89\u529b\u7684\u52b3\u52a8\u529b\u4e5f\u662f\u5fc5\u4e0d\u53ef\u5c11\u7684\u3002\u6b64\u5916\uff0c\u6559\u80b2\u5bf9\u4e8e\u793e\u4f1a\u53d1\u5c55\u4e5f\u81f3\u5173\u91cd\u8981\uff0c\u5b83\u80fd\u5e2e\u52a9\u4e2a\u4eba\u7406\u89e3\u548c\u6b23\u8d4f\u4e0d\u540c\uff0c\u4fc3\u8fdb\u793e\u4f1a\u51dd\u805a\u529b\uff0c\u8425\u9020\u793e\u4f1a\u5171\u540c\u4f53\u610f\u8bc6\u3002\",\n  \"\u6b64\u5916\uff0c\u6559\u80b2\u662f\u6bcf\u4e2a\u4eba\u90fd\u5e94\u8be5\u53ef\u4ee5\u63a5\u53d7\u7684\u57fa\u672c\u6743\u5229\u3002\u5b83\u662f\u516c\u6b63\u548c\u516c\u6b63\u793e\u4f1a\u7684\u57fa\u7840\uff0c\u5728\u5176\u4e2d\u6bcf\u4e2a\u4eba\u90fd\u6709\u5e73\u7b49\u7684\u673a\u4f1a\u53d6\u5f97\u6210\u529f\u3002\u6559\u80b2\u4e3a\u4e2a\u4eba\u63d0\u4f9b\u4e86\u6539\u5584\u751f\u6d3b\u6240\u9700\u7684\u5de5\u5177\u548c\u8d44\u6e90\uff0c\u5b83\u662f\u51cf\u5c11\u8d2b\u56f0\u548c\u4e0d\u5e73\u7b49\u7684\u5f3a\u6709\u529b\u624b\u6bb5\u3002\u6559\u80b2\u5bf9\u4e8e\u4fc3\u8fdb\u6c11\u4e3b\u3001\u4eba\u6743\u548c\u793e\u4f1a\u6b63\u4e49\u4e5f\u81f3\u5173\u91cd\u8981\u3002\",\n  \"\u603b\u4e4b\uff0c\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u81f3\u5173\u91cd\u8981\u7684\u4e00\u4e2a\u65b9\u9762\uff0c\u5bf9\u4e2a\u4eba\u3001\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u5177\u6709\u6df1\u8fdc\u7684\u5f71\u54cd\u3002\u5b83\u662f\u4e00\u79cd\u5de5\u5177\uff0c\u53ef\u4ee5\u4f7f\u4e2a\u4eba\u53d1\u6325\u6f5c\u529b\uff0c\u6539\u5584\u751f\u6d3b\u8d28\u91cf\uff0c\u4e3a\u793e\u4f1a\u505a\u51fa\u8d21\u732e\u3002\u6559\u80b2\u662f\u6bcf\u4e2a\u4eba\u90fd\u5e94\u8be5\u4eab\u6709\u7684\u57fa\u672c\u6743\u5229\uff0c\u4e5f\u662f\u521b\u9020\u516c\u5e73\u6b63\u4e49\u793e\u4f1a\u5fc5\u4e0d\u53ef\u5c11\u7684\u4e00\u90e8\u5206\u3002\u56e0\u6b64\uff0c\u6211\u4eec\u91cd\u89c6\u6559\u80b2\u3001\u786e\u4fdd\u6bcf\u4e2a\u4eba\u90fd\u80fd\u591f\u63a5\u53d7\u4f18\u8d28\u6559\u80b2\u5341\u5206\u91cd\u8981\u3002\",\n];\nconst newTexts = [\n  \"\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u6700\u5173\u952e\u7684\u65b9\u9762\u4e4b\u4e00\u3002\u5b83\u662f\u5e2e\u52a9\u4eba\u4eec\u83b7\u5f97\u77e5\u8bc6\u3001\u6280\u80fd\u3001\u4ef7\u503c\u89c2\u548c\u6001\u5ea6\u3001\u4e3a\u5176\u4e2a\u4eba\u4e0e\u804c\u4e1a\u53d1\u5c55\u6240\u5fc5\u9700\u7684\u5de5\u5177\u3002\u6559\u80b2\u5e76\u4e0d\u4ec5\u9650\u4e8e\u5b66\u6821\u548c\u5927\u5b66\uff1b\u5b83\u53ef\u4ee5\u5728\u4efb\u4f55\u5730\u65b9\u548c\u4efb\u4f55\u65f6\u95f4\u53d1\u751f\u3002\u6559\u80b2\u662f\u4e00\u79cd\u4e0d\u95f4\u65ad\u7684\u8fc7\u7a0b\uff0c\u53ef\u4ee5\u5e2e\u52a9\u4e2a\u4eba\u53d1\u6325\u6f5c\u529b\uff0c\u63d0\u9ad8\u4ed6\u4eec\u7684\u751f\u6d3b\u8d28\u91cf\uff0c\u5e76\u4e3a\u793e\u4f1a\u505a\u51fa\u8d21\u732e\u3002\",\n  \"\u6559\u80b2\u5bf9\u4e2a\u4eba\u7684\u6210\u957f\u548c\u53d1\u5c55\u81f3\u5173\u91cd\u8981\u3002\u5b83\u5e2e\u52a9\u4e2a\u4eba\u66f4\u597d\u5730\u4e86\u89e3\u81ea\u5df1\u7684\u4f18\u70b9\u3001\u7f3a\u70b9\u3001\u5174\u8da3\u548c\u4ef7\u503c\u89c2\u3002\u5b83\u4f7f\u4e2a\u4eba\u80fd\u591f\u8bbe\u5b9a\u76ee\u6807\u3001\u505a\u51fa\u660e\u667a\u7684\u51b3\u7b56\uff0c\u5e76\u5bf9\u81ea\u5df1\u7684\u884c\u4e3a\u8d1f\u8d23\u3002\u6559\u80b2\u8fd8\u80fd\u5e2e\u52a9\u4e2a\u4eba\u53d1\u5c55\u6279\u5224\u6027\u601d\u7ef4\u80fd\u529b\u3001\u521b\u9020\u529b\u548c\u89e3\u51b3\u95ee\u9898\u7684\u80fd\u529b\u3002\u8fd9\u4e9b\u6280\u80fd\u5728\u4eca\u5929\u7684\u4e16\u754c\u4e2d\u81f3\u5173\u91cd\u8981\uff0c\u800c\u5728\u8fd9\u4e2a\u4e16\u754c\u4e2d\uff0c\u9002\u5e94\u53d8\u5316\u548c\u89e3\u51b3\u590d\u6742\u95ee\u9898\u7684\u80fd\u529b\u975e\u5e38\u91cd\u8981\u3002\",\n  \"\u6559\u80b2\u4e5f\u5728\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u4e2d\u626e\u6f14\u4e86\u91cd\u8981\u89d2\u8272\u3002\u5b83\u5e2e\u52a9\u4e2a\u4eba\u83b7\u5f97\u5c31\u4e1a\u548c\u521b\u4e1a\u6240\u5fc5\u987b\u7684\u6280\u80fd\u548c\u77e5\u8bc6\u3002\u6559\u80b2\u5bf9\u4e8e\u521b\u9020\u4e00\u4e2a\u6709\u6280\u80fd\u548c\u80fd\u529b\u7684\u52b3\u52a8\u529b\u6765\u4fc3\u8fdb\u7ecf\u6d4e\u589e\u957f\u548c\u53d1\u5c55\u81f3\u5173\u91cd\u8981\u3002\u5b83\u4e5f\u5bf9\u4e8e\u793e\u4f1a\u53d1\u5c55\u81f3\u5173\u91cd\u8981\uff0c\u56e0\u4e3a\u5b83\u5e2e\u52a9\u4e2a\u4eba\u7406\u89e3\u548c\u6b23\u8d4f\u591a\u6837\u6027\uff0c\u4fc3\u8fdb\u793e\u4f1a\u51dd\u805a\u529b\uff0c\u5e76\u57f9\u517b\u793e\u533a\u610f\u8bc6\u3002\",\n  \"\u6b64\u5916\uff0c\u6559\u80b2\u662f\u4e00\u4e2a\u5e94\u8be5\u4e3a\u6bcf\u4e2a\u4eba\u6240\u63a5\u89e6\u5230\u7684\u57fa\u672c\u6743\u5229\u3002\u5b83\u662f\u4e00\u4e2a\u516c\u5e73\u548c\u6b63\u4e49\u793e\u4f1a\u7684\u57fa\u7840\uff0c\u8ba9\u6bcf\u4e2a\u4eba\u90fd\u6709\u5e73\u7b49\u7684\u673a\u4f1a\u53bb\u83b7\u5f97\u6210\u529f\u3002\u6559\u80b2\u4e3a\u4e2a\u4eba\u63d0\u4f9b\u4e86\u6539\u5584\u751f\u6d3b\u6240\u9700\u7684\u5de5\u5177\u548c\u8d44\u6e90\uff0c\u5b83\u662f\u51cf\u5c11\u8d2b\u56f0\u548c\u4e0d\u5e73\u7b49\u7684\u6709\u529b\u624b\u6bb5\u3002\u6559\u80b2\u5bf9\u4e8e\u4fc3\u8fdb\u6c11\u4e3b\u3001\u4eba\u6743\u548c\u793e\u4f1a\u516c\u6b63\u4e5f\u81f3\u5173\u91cd\u8981\u3002\",\n  \"\u603b\u4e4b\uff0c\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u81f3\u5173\u91cd\u8981\u7684\u4e00\u65b9\u9762\uff0c\u5bf9\u4e2a\u4eba\u3001\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u5177\u6709\u6df1\u8fdc\u7684\u5f71\u54cd\u3002\u5b83\u662f\u4e00\u79cd\u5de5\u5177\uff0c\u4f7f\u4e2a\u4eba\u80fd\u591f\u53d1\u5c55\u4ed6\u4eec\u7684\u6f5c\u529b\uff0c\u6539\u5584\u4ed6\u4eec\u7684\u751f\u6d3b\u8d28\u91cf\u5e76\u4e3a\u793e\u4f1a\u505a\u51fa\u8d21\u732e\u3002\u6559\u80b2\u662f\u4e00\u9879\u57fa\u672c\u6743\u5229\uff0c\u5e94\u8be5\u4e3a\u6bcf\u4e2a\u4eba\u6240\u80fd\u63a5\u89e6\u5230\uff0c\u5b83\u5bf9\u4e8e\u521b\u9020\u4e00\u4e2a\u516c\u6b63\u548c\u516c\u5e73\u7684\u793e\u4f1a\u662f\u5fc5\u8981\u7684\u3002\u56e0\u6b64\uff0c\u6211\u4eec\u6295\u8d44\u4e8e\u6559\u80b2\u5e76\u786e\u4fdd\u6bcf\u4e2a\u4eba\u90fd\u6709\u63a5\u53d7\u4f18\u8d28\u6559\u80b2\u7684\u673a\u4f1a\u662f\u81f3\u5173\u91cd\u8981\u7684\u3002\",\n];\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nfunction matchIndex(text) {\n  let idx = oldTexts.indexOf(text);\n  if (idx !== -1) return idx;\n  const trimmed = text.trim();\n  for (let i = 0; i < oldTexts.length; i++) {\n    if (oldTexts[i].trim() === trimmed) return i;\n  }\n  return -1;\n}\n\nlet replaced = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const idx = matchIndex(para.text);\n  if (idx !== -1) {\n    para.insertText(newTexts[idx], \"Replace\");\n    replaced++;\n  }\n}\nawait context.sync();\n\nreturn \"replaced=\" + replaced;\n", "ps1": "$oldTexts = @(\n  \"\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u6700\u5173\u952e\u7684\u65b9\u9762\u4e4b\u4e00\u3002 \u5b83\u662f\u4e00\u4e2a\u5e2e\u52a9\u4eba\u4eec\u83b7\u5f97\u77e5\u8bc6\u3001\u6280\u80fd\u3001\u4ef7\u503c\u89c2\u548c\u6001\u5ea6\u7684\u5de5\u5177\uff0c\u8fd9\u4e9b\u5bf9\u4e8e\u4ed6\u4eec\u7684\u4e2a\u4eba\u548c\u804c\u4e1a\u53d1\u5c55\u662f\u5fc5\u8981\u7684\u3002\u6559\u80b2\u5e76\u4e0d\u4ec5\u9650\u4e8e\u5b66\u6821\u548c\u5927\u5b66\uff1b\u5b83\u53ef\u4ee5\u5728\u4efb\u4f55\u5730\u65b9\u548c\u4efb\u4f55\u65f6\u5019\u53d1\u751f\u3002\u6559\u80b2\u662f\u4e00\u4e2a\u6301\u7eed\u7684\u8fc7\u7a0b\uff0c\u5e2e\u52a9\u4e2a\u4eba\u53d1\u6325\u6f5c\u529b\uff0c\u6539\u5584\u751f\u6d3b\u8d28\u91cf\uff0c\u5e76\u4e3a\u793e\u4f1a\u4f5c\u51fa\u8d21\u732e\u3002\",\n  \"\u6559\u80b2\u5bf9\u4e2a\u4eba\u6210\u957f\u548c\u53d1\u5c55\u81f3\u5173\u91cd\u8981\u3002\u5b83\u5e2e\u52a9\u4e2a\u4eba\u66f4\u597d\u5730\u4e86\u89e3\u81ea\u5df1\u7684\u4f18\u70b9\u3001\u7f3a\u70b9\u3001\u5174\u8da3\u548c\u4ef7\u503c\u89c2\uff0c\u4f7f\u4ed6\u4eec\u80fd\u591f\u5236\u5b9a\u76ee\u6807\u3001\u505a\u51fa\u660e\u667a\u51b3\u7b56\u5e76\u5bf9\u81ea\u5df1\u7684\u884c\u52a8\u8d1f\u8d23\u3002\u6559\u80b2\u8fd8\u5e2e\u52a9\u4e2a\u4eba\u57f9\u517b\u6279\u5224\u6027\u601d\u7ef4\u80fd\u529b\u3001\u521b\u9020\u529b\u548c\u89e3\u51b3\u95ee\u9898\u7684\u80fd\u529b\u3002\u8fd9\u4e9b\u6280\u80fd\u5bf9\u4e8e\u5728\u4eca\u5929\u7684\u4e16\u754c\u4e2d\u6210\u529f\u81f3\u5173\u91cd\u8981\uff0c\u8fd9\u4e2a\u4e16\u754c\u4e2d\uff0c\u9002\u5e94\u53d8\u5316\u548c\u89e3\u51b3\u590d\u6742\u95ee\u9898\u7684\u80fd\u529b\u662f\u81f3\u5173\u91cd\u8981\u7684\u3002\",\n  \"\u6559\u80b2\u5728\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u4e2d\u4e5f\u626e\u6f14\u7740\u91cd\u8981\u89d2\u8272\u3002\u5b83\u80fd\u5e2e\u52a9\u4e2a\u4eba\u83b7\u53d6\u5c31\u4e1a\u548c\u521b\u4e1a\u6240\u5fc5\u9700\u7684\u6280\u80fd\u548c\u77e5\u8bc6\u3002\u6559\u80b2\u5bf9\u4e8e\u57f9\u517b\u80fd\u591f\u4e3a\u7ecf\u6d4e\u589e\u957f\u548c\u53d1\u5c55\u8d21\u732e\u7684\u6709\u6280\u80fd\u3001\u6709\u7ade\u4e89\u529b\u7684\u52b3\u52a8\u529b\u4e5f\u662f\u5fc5\u4e0d\u53ef\u5c11\u7684\u3002\u6b64\u5916\uff0c\u6559\u80b2\u5bf9\u4e8e\u793e\u4f1a\u53d1\u5c55\u4e5f\u81f3\u5173\u91cd\u8981\uff0c\u5b83\u80fd\u5e2e\u52a9\u4e2a\u4eba\u7406\u89e3\u548c\u6b23\u8d4f\u4e0d\u540c\uff0c\u4fc3\u8fdb\u793e\u4f1a\u51dd\u805a\u529b\uff0c\u8425\u9020\u793e\u4f1a\u5171\u540c\u4f53\u610f\u8bc6\u3002\",\n  \"\u6b64\u5916\uff0c\u6559\u80b2\u662f\u6bcf\u4e2a\u4eba\u90fd\u5e94\u8be5\u53ef\u4ee5\u63a5\u53d7\u7684\u57fa\u672c\u6743\u5229\u3002\u5b83\u662f\u516c\u6b63\u548c\u516c\u6b63\u793e\u4f1a\u7684\u57fa\u7840\uff0c\u5728\u5176\u4e2d\u6bcf\u4e2a\u4eba\u90fd\u6709\u5e73\u7b49\u7684\u673a\u4f1a\u53d6\u5f97\u6210\u529f\u3002\u6559\u80b2\u4e3a\u4e2a\u4eba\u63d0\u4f9b\u4e86\u6539\u5584\u751f\u6d3b\u6240\u9700\u7684\u5de5\u5177\u548c\u8d44\u6e90\uff0c\u5b83\u662f\u51cf\u5c11\u8d2b\u56f0\u548c\u4e0d\u5e73\u7b49\u7684\u5f3a\u6709\u529b\u624b\u6bb5\u3002\u6559\u80b2\u5bf9\u4e8e\u4fc3\u8fdb\u6c11\u4e3b\u3001\u4eba\u6743\u548c\u793e\u4f1a\u6b63\u4e49\u4e5f\u81f3\u5173\u91cd\u8981\u3002\",\n  \"\u603b\u4e4b\uff0c\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u81f3\u5173\u91cd\u8981\u7684\u4e00\u4e2a\u65b9\u9762\uff0c\u5bf9\u4e2a\u4eba\u3001\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u5177\u6709\u6df1\u8fdc\u7684\u5f71\u54cd\u3002\u5b83\u662f\u4e00\u79cd\u5de5\u5177\uff0c\u53ef\u4ee5\u4f7f\u4e2a\u4eba\u53d1\u6325\u6f5c\u529b\uff0c\u6539\u5584\u751f\u6d3b\u8d28\u91cf\uff0c\u4e3a\u793e\u4f1a\u505a\u51fa\u8d21\u732e\u3002\u6559\u80b2\u662f\u6bcf\u4e2a\u4eba\u90fd\u5e94\u8be5\u4eab\u6709\u7684\u57fa\u672c\u6743\u5229\uff0c\u4e5f\u662f\u521b\u9020\u516c\u5e73\u6b63\u4e49\u793e\u4f1a\u5fc5\u4e0d\u53ef\u5c11\u7684\u4e00\u90e8\u5206\u3002\u56e0\u6b64\uff0c\u6211\u4eec\u91cd\u89c6\u6559\u80b2\u3001\u786e\u4fdd\u6bcf\u4e2a\u4eba\u90fd\u80fd\u591f\u63a5\u53d7\u4f18\u8d28\u6559\u80b2\u5341\u5206\u91cd\u8981\u3002\",\n)\n$newTexts = @(\n  \"\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u6700\u5173\u952e\u7684\u65b9\u9762\u4e4b\u4e00\u3002\u5b83\u662f\u5e2e\u52a9\u4eba\u4eec\u83b7\u5f97\u77e5\u8bc6\u3001\u6280\u80fd\u3001\u4ef7\u503c\u89c2\u548c\u6001\u5ea6\u3001\u4e3a\u5176\u4e2a\u4eba\u4e0e\u804c\u4e1a\u53d1\u5c55\u6240\u5fc5\u9700\u7684\u5de5\u5177\u3002\u6559\u80b2\u5e76\u4e0d\u4ec5\u9650\u4e8e\u5b66\u6821\u548c\u5927\u5b66\uff1b\u5b83\u53ef\u4ee5\u5728\u4efb\u4f55\u5730\u65b9\u548c\u4efb\u4f55\u65f6\u95f4\u53d1\u751f\u3002\u6559\u80b2\u662f\u4e00\u79cd\u4e0d\u95f4\u65ad\u7684\u8fc7\u7a0b\uff0c\u53ef\u4ee5\u5e2e\u52a9\u4e2a\u4eba\u53d1\u6325\u6f5c\u529b\uff0c\u63d0\u9ad8\u4ed6\u4eec\u7684\u751f\u6d3b\u8d28\u91cf\uff0c\u5e76\u4e3a\u793e\u4f1a\u505a\u51fa\u8d21\u732e\u3002\",\n  \"\u6559\u80b2\u5bf9\u4e2a\u4eba\u7684\u6210\u957f\u548c\u53d1\u5c55\u81f3\u5173\u91cd\u8981\u3002\u5b83\u5e2e\u52a9\u4e2a\u4eba\u66f4\u597d\u5730\u4e86\u89e3\u81ea\u5df1\u7684\u4f18\u70b9\u3001\u7f3a\u70b9\u3001\u5174\u8da3\u548c\u4ef7\u503c\u89c2\u3002\u5b83\u4f7f\u4e2a\u4eba\u80fd\u591f\u8bbe\u5b9a\u76ee\u6807\u3001\u505a\u51fa\u660e\u667a\u7684\u51b3\u7b56\uff0c\u5e76\u5bf9\u81ea\u5df1\u7684\u884c\u4e3a\u8d1f\u8d23\u3002\u6559\u80b2\u8fd8\u80fd\u5e2e\u52a9\u4e2a\u4eba\u53d1\u5c55\u6279\u5224\u6027\u601d\u7ef4\u80fd\u529b\u3001\u521b\u9020\u529b\u548c\u89e3\u51b3\u95ee\u9898\u7684\u80fd\u529b\u3002\u8fd9\u4e9b\u6280\u80fd\u5728\u4eca\u5929\u7684\u4e16\u754c\u4e2d\u81f3\u5173\u91cd\u8981\uff0c\u800c\u5728\u8fd9\u4e2a\u4e16\u754c\u4e2d\uff0c\u9002\u5e94\u53d8\u5316\u548c\u89e3\u51b3\u590d\u6742\u95ee\u9898\u7684\u80fd\u529b\u975e\u5e38\u91cd\u8981\u3002\",\n  \"\u6559\u80b2\u4e5f\u5728\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u4e2d\u626e\u6f14\u4e86\u91cd\u8981\u89d2\u8272\u3002\u5b83\u5e2e\u52a9\u4e2a\u4eba\u83b7\u5f97\u5c31\u4e1a\u548c\u521b\u4e1a\u6240\u5fc5\u987b\u7684\u6280\u80fd\u548c\u77e5\u8bc6\u3002\u6559\u80b2\u5bf9\u4e8e\u521b\u9020\u4e00\u4e2a\u6709\u6280\u80fd\u548c\u80fd\u529b\u7684\u52b3\u52a8\u529b\u6765\u4fc3\u8fdb\u7ecf\u6d4e\u589e\u957f\u548c\u53d1\u5c55\u81f3\u5173\u91cd\u8981\u3002\u5b83\u4e5f\u5bf9\u4e8e\u793e\u4f1a\u53d1\u5c55\u81f3\u5173\u91cd\u8981\uff0c\u56e0\u4e3a\u5b83\u5e2e\u52a9\u4e2a\u4eba\u7406\u89e3\u548c\u6b23\u8d4f\u591a\u6837\u6027\uff0c\u4fc3\u8fdb\u793e\u4f1a\u51dd\u805a\u529b\uff0c\u5e76\u57f9\u517b\u793e\u533a\u610f\u8bc6\u3002\",\n  \"\u6b64\u5916\uff0c\u6559\u80b2\u662f\u4e00\u4e2a\u5e94\u8be5\u4e3a\u6bcf\u4e2a\u4eba\u6240\u63a5\u89e6\u5230\u7684\u57fa\u672c\u6743\u5229\u3002\u5b83\u662f\u4e00\u4e2a\u516c\u5e73\u548c\u6b63\u4e49\u793e\u4f1a\u7684\u57fa\u7840\uff0c\u8ba9\u6bcf\u4e2a\u4eba\u90fd\u6709\u5e73\u7b49\u7684\u673a\u4f1a\u53bb\u83b7\u5f97\u6210\u529f\u3002\u6559\u80b2\u4e3a\u4e2a\u4eba\u63d0\u4f9b\u4e86\u6539\u5584\u751f\u6d3b\u6240\u9700\u7684\u5de5\u5177\u548c\u8d44\u6e90\uff0c\u5b83\u662f\u51cf\u5c11\u8d2b\u56f0\u548c\u4e0d\u5e73\u7b49\u7684\u6709\u529b\u624b\u6bb5\u3002\u6559\u80b2\u5bf9\u4e8e\u4fc3\u8fdb\u6c11\u4e3b\u3001\u4eba\u6743\u548c\u793e\u4f1a\u516c\u6b63\u4e5f\u81f3\u5173\u91cd\u8981\u3002\",\n  \"\u603b\u4e4b\uff0c\u6559\u80b2\u662f\u4eba\u7c7b\u751f\u6d3b\u4e2d\u81f3\u5173\u91cd\u8981\u7684\u4e00\u65b9\u9762\uff0c\u5bf9\u4e2a\u4eba\u3001\u793e\u4f1a\u548c\u7ecf\u6d4e\u53d1\u5c55\u5177\u6709\u6df1\u8fdc\u7684\u5f71\u54cd\u3002\u5b83\u662f\u4e00\u79cd\u5de5\u5177\uff0c\u4f7f\u4e2a\u4eba\u80fd\u591f\u53d1\u5c55\u4ed6\u4eec\u7684\u6f5c\u529b\uff0c\u6539\u5584\u4ed6\u4eec\u7684\u751f\u6d3b\u8d28\u91cf\u5e76\u4e3a\u793e\u4f1a\u505a\u51fa\u8d21\u732e\u3002\u6559\u80b2\u662f\u4e00\u9879\u57fa\u672c\u6743\u5229\uff0c\u5e94\u8be5\u4e3a\u6bcf\u4e2a\u4eba\u6240\u80fd\u63a5\u89e6\u5230\uff0c\u5b83\u5bf9\u4e8e\u521b\u9020\u4e00\u4e2a\u516c\u6b63\u548c\u516c\u5e73\u7684\u793e\u4f1a\u662f\u5fc5\u8981\u7684\u3002\u56e0\u6b64\uff0c\u6211\u4eec\u6295\u8d44\u4e8e\u6559\u80b2\u5e76\u786e\u4fdd\u6bcf\u4e2a\u4eba\u90fd\u6709\u63a5\u53d7\u4f18\u8d28\u6559\u80b2\u7684\u673a\u4f1a\u662f\u81f3\u5173\u91cd\u8981\u7684\u3002\",\n)\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replaced = 0\nfor ($i = 0; $i -lt $oldTexts.Length; $i++) {\n    $old = $oldTexts[$i]\n    $new = $newTexts[$i]\n\n    $range = $d.Content\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n    if ($found) {\n        $replaced++\n    } else {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n\nWrite-Output \"replaced=$replaced\"\n"}
